$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-35 down to 25-36
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with data
$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 45135
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100112013
$ws.Cells.Item(24, 7).Value = "Alcachofa"
$ws.Cells.Item(24, 8).Value = "Madrigal"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 500
$ws.Cells.Item(24, 11).Value = 10000
$ws.Cells.Item(24, 12).Value = 11000
$ws.Cells.Item(24, 13).Value = 10500
$ws.Cells.Item(24, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(24, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(24, 16).Value = 262
$ws.Cells.Item(24, 17).Value = 40
$ws.Cells.Item(24, 18).Value = "Hortaliza"
